# Auto-applies the OOXML-level edit described by the commit diff.
#
# The whole body (minus the trailing <w:sectPr>, which Word's
# Range "Content" never includes and therefore doesn't need to be
# restated -- the original section properties are preserved
# automatically) is replaced in one shot via Range.InsertXML,
# passing an explicit WordprocessingML package-part wrapper. This
# lets every run/paragraph split, proofErr span removal, bookmark
# relocation, list numbering, and highlight/bold run exactly match
# the target OOXML instead of approximating it through many small
# Find/Replace calls.
#
# Word's COM model never lets a Range delete the document's very
# last paragraph mark, so InsertXML always leaves one trailing
# paragraph behind no matter what XML is supplied for it. To work
# around that we give the final (empty) target paragraph a
# recognizable placeholder run, insert everything, then clear that
# placeholder's text with an ordinary Find/Replace so the last
# paragraph ends up genuinely empty, matching the target.

$d = $word.ActiveDocument
$placeholder = "ZZZ_TRAILING_PLACEHOLDER_ZZZ"

$target = $d.Content
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 wp14"><w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Summary of my error checking </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">File last updated </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> DATE \@ "d/MM/yyyy h:mm:ss am/pm" </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>9/10/2014 3:32:43 PM</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p><w:p><w:r><w:t>Re</w:t></w:r><w:r><w:t xml:space="preserve">curring </w:t></w:r><w:r><w:t>problems that I think deriv</w:t></w:r><w:r><w:t>e from problems with plant maps or the program:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>BAER</w:t></w:r><w:r><w:t xml:space="preserve">, transition from </w:t></w:r><w:r><w:t>cone_green_01</w:t></w:r><w:r><w:t xml:space="preserve"> to cone_brown has problems in some plants (see EHW_notes_FinDev)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>COER</w:t></w:r><w:r><w:t xml:space="preserve">, nothing in the </w:t></w:r><w:r><w:t>investment table ends with "inflorescence_stalk_in_fruit"</w:t></w:r><w:r><w:t>, yet this should exist for most plants</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">GRBU, </w:t></w:r><w:r><w:t>the category “bud_tiny” needs to be added to the start of the GRBU plant map; this has caused lots of “errors” (so to speak)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">HEPU, </w:t></w:r><w:r><w:t>there are “bud_small” correctly listed in the FinDev file, but none in the lost parts file, even though the ones in the FinDev file are indeed lost parts.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>PUTU</w:t></w:r><w:r><w:t xml:space="preserve">, various problems with transition from finished_flower_stigma to fruit_large_immature and onto seed/seed_pod, resulting in </w:t></w:r><w:r><w:t xml:space="preserve">"fruit_large_immature" </w:t></w:r><w:r><w:t xml:space="preserve">incorrectly appearing </w:t></w:r><w:r><w:t xml:space="preserve">as a </w:t></w:r><w:r><w:t xml:space="preserve">FinDev </w:t></w:r><w:r><w:t>part</w:t></w:r><w:r><w:t xml:space="preserve"> and lots of “</w:t></w:r><w:r><w:t>errors</w:t></w:r><w:r><w:t xml:space="preserve"> in many "finished_flower_stigma” entries in error file</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:r><w:t>Other questions/problems:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>With PELA and BAER, where young fruit with dimensions pass through many stages as they grow, the program seems to not always “pick out” the correct predecessor. See examples in EHW_notes_INV. I don’t have my head around a fix, but for these species I can easily re-enter some numbers differently if it would be helpful.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Only other recurring problem has to do with the fixed splitting of resources from 1 part into 2 (or more) subsequent parts, resulting in negative investment. Lots of cases for BOLE, PUTU</w:t></w:r><w:r><w:t>, a few for various other species</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>With PHPH, flower petals small are in the map as an intermediate stage, but actually represent the endpoint for the individuals where I use the category. Does it matter for any other calculations that those plant parts are considered “lost” instead of a true endpoint?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>For a handful of species, some part that should hasn’t showed up in either FinDev or as Lost. While I’ve looked at every part and made sure it was present and its approximate predecessor’s made sense to me, I haven’t actually added up all the numbers to make sure the correct number of other parts flow on to their final location, Do you think this would be a good idea, given the “missing” parts – or is it likely I’ve caught the few little mistakes? (For instance, with GRSP, there is no FinDev for Fruit_young – doesn’t surprise me – but there are 28 counts in the “To” column, but only 56 in the “From” column. This is of course a multiple of 2, so maybe different sort of error.) But I’m just wondering if I – perhaps simplified by a computer computer script – should make sure the different between the “From” and “To” counts is the number in “FinDev”</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Species thoroughly checked: </w:t></w:r><w:r><w:t xml:space="preserve">BAER, BOLE, COER, EPMI, </w:t></w:r><w:r><w:t xml:space="preserve">GRSP, </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">HATE, </w:t></w:r><w:r><w:t>HEPU, LEES</w:t></w:r><w:r><w:t xml:space="preserve">, PELA, </w:t></w:r><w:r><w:t xml:space="preserve">PEPU, PHPH, </w:t></w:r><w:r><w:t>PILI, PUTU</w:t></w:r><w:r><w:t>; GRBU not completely checked at this point</w:t></w:r></w:p><w:p><w:r><w:t>Species for Lizzy to recheck in the future:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>PEPU – after rest of flowering parts collected</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>PUTU</w:t></w:r><w:r><w:t>, GRBU</w:t></w:r><w:r><w:t xml:space="preserve"> – after plant map fixed</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Still for Lizzy to fix:</w:t></w:r></w:p><w:p><w:r><w:t>-adjust some fruit_immature sizes for GRSP to make sure progression is always positive</w:t></w:r></w:p><w:p><w:r><w:t>-GRBU infl_bud is still smaller than the stalk; I still need to look into this</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>And I can’t remember why certain plants have errors and fail to run at all…</w:t></w:r></w:p><w:p><w:r><w:t>ZZZ_TRAILING_PLACEHOLDER_ZZZ</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)

$cleanup = $d.Content
$found = $cleanup.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $cleanup.Text = ""
} else {
    throw "edit.ps1: failed to locate trailing placeholder run after InsertXML"
}

Write-Host "edit applied"
